$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 6 new rows for the additional worker (TOMAS IGNACIO ZAMBRANO BOLIVAR),
# pushing the signature block (old rows 39:40) down to 45:46.
$ws.Rows("35:40").Insert()

# Restore proper table-body formatting for the newly-inserted rows: the
# regular mid-table row style for 35:39, and the closing (bottom-border) row
# style -- which used to belong to the old last row (34) -- for row 40.
$ws.Range("B16:J16").Copy()
$ws.Range("B35:J39").PasteSpecial(-4122)
$ws.Range("B34:J34").Copy()
$ws.Range("B40:J40").PasteSpecial(-4122)

# Updated totals
$ws.Range("E11").Value = 723390
$ws.Range("C13").Value = 5

# New account-statement body: 5 workers x their period rows.
$data = @(
  @("CC", "1047454838", "MELISSA ANDREA COGOLLO MOLINA", "1906", 17708, 781242),
  @("CC", "1047454838", "MELISSA ANDREA COGOLLO MOLINA", "1905", 31249, 781242),
  @("CC", "1047454838", "MELISSA ANDREA COGOLLO MOLINA", "1904", 31249, 781242),
  @("CC", "1047454838", "MELISSA ANDREA COGOLLO MOLINA", "1903", 31249, 781242),
  @("CC", "1047454838", "MELISSA ANDREA COGOLLO MOLINA", "1902", 31249, 781242),
  @("CC", "1047454838", "MELISSA ANDREA COGOLLO MOLINA", "1901", 31249, 781242),
  @("CC", "45534413",   "LILA MARINA ZAMBRANO BOLIVAR",  "1609", 27578, 737717),
  @("CC", "1047365598", "GUILLERMO JOSE ZAMBRANO BOLIVAR","1906", 17708, 781242),
  @("CC", "1047365598", "GUILLERMO JOSE ZAMBRANO BOLIVAR","1905", 31249, 781242),
  @("CC", "1047365598", "GUILLERMO JOSE ZAMBRANO BOLIVAR","1904", 31249, 781242),
  @("CC", "1047365598", "GUILLERMO JOSE ZAMBRANO BOLIVAR","1903", 31249, 781242),
  @("CC", "1047365598", "GUILLERMO JOSE ZAMBRANO BOLIVAR","1902", 31249, 781242),
  @("CC", "1047365598", "GUILLERMO JOSE ZAMBRANO BOLIVAR","1901", 31249, 781242),
  @("CC", "1047365597", "TOMAS IGNACIO ZAMBRANO BOLIVAR", "1906", 17708, 781242),
  @("CC", "1047365597", "TOMAS IGNACIO ZAMBRANO BOLIVAR", "1905", 31249, 781242),
  @("CC", "1047365597", "TOMAS IGNACIO ZAMBRANO BOLIVAR", "1904", 31249, 781242),
  @("CC", "1047365597", "TOMAS IGNACIO ZAMBRANO BOLIVAR", "1903", 31249, 781242),
  @("CC", "1047365597", "TOMAS IGNACIO ZAMBRANO BOLIVAR", "1902", 31249, 781242),
  @("CC", "1047365597", "TOMAS IGNACIO ZAMBRANO BOLIVAR", "1901", 31249, 781242),
  @("CC", "1047432113", "VICTOR ANDRES COGOLLO MOLINA",   "1906", 17708, 781242),
  @("CC", "1047432113", "VICTOR ANDRES COGOLLO MOLINA",   "1905", 31249, 781242),
  @("CC", "1047432113", "VICTOR ANDRES COGOLLO MOLINA",   "1904", 31249, 781242),
  @("CC", "1047432113", "VICTOR ANDRES COGOLLO MOLINA",   "1903", 31249, 781242),
  @("CC", "1047432113", "VICTOR ANDRES COGOLLO MOLINA",   "1902", 31249, 781242),
  @("CC", "1047432113", "VICTOR ANDRES COGOLLO MOLINA",   "1901", 31249, 781242)
)

$startRow = 16
for ($i = 0; $i -lt $data.Count; $i++) {
  $r = $startRow + $i
  $row = $data[$i]
  $ws.Range("B$r").Value = $row[0]
  $ws.Range("C$r").Value = $row[1]
  $ws.Range("D$r").Value = $row[2]
  $ws.Range("E$r").Value = $row[3]
  $ws.Range("F$r").Value = $row[4]
  $ws.Range("G$r").Value = $row[5]
}

Write-Host "done"
